$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column-B labels for rows 4-29 end up shifted down by two positions versus
# the prior run (two new entries - "Holden" and "Rizzie Spiral" - land at
# rows 4-5, pushing every following label down two rows) and "Thomas Hex" is
# renamed to "Matthies Hex". Row/Index (col A) and the per-column simulation
# results (cols C:T) are left as they were.
$labels = @(
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 4, 2).Value = $labels[$i]
}

# Append two new simulation rows (28: Michael-CCHex, 29: Michael-SNHex),
# reusing row 29's formatting for the new index cells in column A.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("C30:T30").Value = 1

$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Michael-SNHex"
$ws.Range("C31:T31").Value = 1
